# Epoch Accuracy.xlsx -- rerun notebook output refresh (M07 Froze Encoder 123)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Refreshed per-epoch accuracy values in column B (A holds the epoch index / repr text)
$values = @{
    2 = 0.9375
    3 = 0.890625
    4 = 0.890625
    5 = 0.859375
    6 = 0.828125
    7 = 0.8125
    8 = 0.796875
    9 = 0.71875
    10 = 0.765625
    11 = 0.765625
    12 = 0.765625
    13 = 0.65625
    14 = 0.703125
    15 = 0.640625
    16 = 0.625
    17 = 0.625
    18 = 0.546875
    19 = 0.53125
    20 = 0.515625
    21 = 0.5625
    22 = 0.546875
    23 = 0.5625
    24 = 0.5625
    25 = 0.5625
    26 = 0.5625
    27 = 0.5625
    28 = 0.578125
    29 = 0.578125
    30 = 0.578125
    31 = 0.5625
    32 = 0.5625
    33 = 0.578125
    34 = 0.578125
    35 = 0.5625
    36 = 0.5625
    37 = 0.5625
    38 = 0.5625
    39 = 0.5625
    40 = 0.5625
    41 = 0.5625
    42 = 0.5625
    43 = 0.5625
    44 = 0.5625
    45 = 0.5625
    46 = 0.5625
    47 = 0.5625
    48 = 0.5625
    49 = 0.5625
    50 = 0.5625
    51 = 0.5625
    52 = 0.5625
    53 = 0.5625
    54 = 0.5625
    55 = 0.5625
    56 = 0.5625
    57 = 0.5625
    58 = 0.5625
    59 = 0.5625
    60 = 0.5625
    61 = 0.5625
    62 = 0.5625
    63 = 0.5625
    64 = 0.5625
    65 = 0.5625
    66 = 0.578125
    67 = 0.578125
    68 = 0.578125
    69 = 0.578125
    70 = 0.59375
    71 = 0.59375
    72 = 0.59375
    73 = 0.59375
    74 = 0.59375
    75 = 0.59375
    76 = 0.59375
    77 = 0.59375
    78 = 0.59375
    79 = 0.59375
    80 = 0.59375
    81 = 0.59375
    82 = 0.59375
    83 = 0.59375
    84 = 0.59375
    85 = 0.59375
    86 = 0.59375
    87 = 0.59375
    88 = 0.59375
    89 = 0.59375
    90 = 0.59375
    91 = 0.59375
    92 = 0.59375
    93 = 0.59375
    94 = 0.59375
    95 = 0.59375
    96 = 0.59375
    97 = 0.59375
    98 = 0.59375
    99 = 0.59375
    100 = 0.59375
    101 = 0.59375
    102 = 0.59375
    103 = 0.6875
    104 = 0.71875
    105 = 0.546875
    106 = 0.53125
    107 = 0.609375
    108 = 0.671875
    109 = 0.671875
    110 = 0.578125
    111 = 0.5625
    112 = 0.578125
    113 = 0.640625
    114 = 0.59375
    115 = 0.609375
    116 = 0.625
    117 = 0.5625
    118 = 0.5081967213114754
}
foreach ($r in $values.Keys) {
    $ws.Cells.Item($r, 2).Value = $values[$r]
}

# Rows 102-118 held a literal repr() of the DisplayOutputs object in column A;
# re-running the notebook reallocated it at a new address, so the text changed too.
for ($r = 102; $r -le 118; $r++) {
    $ws.Cells.Item($r, 1).Value = "<__main__.DisplayOutputs object at 0x7fdcf0507e50>"
}

# The sheet was left with everything selected (e.g. Ctrl+A) before saving.
$ws.Cells.Select()
